$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-RowText($rowIndex, $newText) {
    $cell = $t.Rows.Item($rowIndex).Cells.Item(1)
    $cell.Range.Text = $newText
}

# Simple single-run value replacements
Set-RowText 1 "0M"
Set-RowText 2 "0M"
Set-RowText 3 "0M"
Set-RowText 4 "455"
Set-RowText 6 "0.04323"
Set-RowText 7 "0.00129"
Set-RowText 8 "0.00148"
Set-RowText 9 "0.00311"
Set-RowText 10 "0.00758"
Set-RowText 11 "0.01243"
Set-RowText 12 "0.12506"

# Rows whose multi-run/tab content collapses into a single simple value
Set-RowText 44 "99.88"
Set-RowText 45 "0.13"
Set-RowText 46 "108"
